$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$sh.TextFrame.TextRange.Text = "30/05/2021"
Write-Host "done"
